$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 993.2
$ws.Cells.Item(19, 9).Value = 883.75
$ws.Cells.Item(19, 10).Value = 1066.1666
$ws.Cells.Item(19, 11).Value = 883.75
$ws.Cells.Item(19, 12).Value = 1066.1666
$ws.Cells.Item(19, 13).Value = -708.75
$ws.Cells.Item(19, 14).Value = -1416.1666

$ws.Cells.Item(107, 8).Value = 785.1667
$ws.Cells.Item(107, 9).Value = 785.1667
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 785.1667
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 1134.8333
$ws.Cells.Item(107, 14).ClearContents()

$ws.Cells.Item(121, 8).Value = 3744.9546
$ws.Cells.Item(121, 10).Value = 3866.3809
$ws.Cells.Item(121, 12).Value = 11599.1427
$ws.Cells.Item(121, 14).Value = -15093.1427

$ws.Cells.Item(129, 8).Value = 1748.8334
$ws.Cells.Item(129, 9).Value = 1365.6666
$ws.Cells.Item(129, 10).Value = 2132
$ws.Cells.Item(129, 11).Value = 4096.9998
$ws.Cells.Item(129, 12).Value = 6396
$ws.Cells.Item(129, 13).Value = 903.0002000000004
$ws.Cells.Item(129, 14).Value = -16396

$ws.Cells.Item(132, 8).Value = 107847.7
$ws.Cells.Item(132, 9).Value = 247248.73
$ws.Cells.Item(132, 11).Value = 741746.1900000001
$ws.Cells.Item(132, 13).Value = -739216.1900000001

$ws.Cells.Item(137, 8).Value = 3084.5186
$ws.Cells.Item(137, 9).Value = 2560.889
$ws.Cells.Item(137, 10).Value = 3346.3333
$ws.Cells.Item(137, 11).Value = 7682.667
$ws.Cells.Item(137, 12).Value = 10038.9999
$ws.Cells.Item(137, 13).Value = -5132.667
$ws.Cells.Item(137, 14).Value = -15138.9999

$ws.Cells.Item(138, 8).Value = 8071.5386
$ws.Cells.Item(138, 9).Value = 2083
$ws.Cells.Item(138, 10).Value = 9568.673000000001
$ws.Cells.Item(138, 11).Value = 6249
$ws.Cells.Item(138, 12).Value = 28706.019
$ws.Cells.Item(138, 13).Value = -1109
$ws.Cells.Item(138, 14).Value = -38986.019

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).ClearContents()

$ws.Cells.Item(32, 8).Value = 3063.037
$ws.Cells.Item(32, 9).Value = 1715.6666
$ws.Cells.Item(32, 10).Value = 13842
$ws.Cells.Item(32, 11).Value = 1715.6666
$ws.Cells.Item(32, 12).Value = 13842
$ws.Cells.Item(32, 13).Value = -1428.6666
$ws.Cells.Item(32, 14).Value = -14416

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 422.36365
$ws.Cells.Item(80, 9).Value = 527
$ws.Cells.Item(80, 10).Value = 383.125
$ws.Cells.Item(80, 11).Value = 527
$ws.Cells.Item(80, 12).Value = 383.125
$ws.Cells.Item(80, 13).Value = 471
$ws.Cells.Item(80, 14).Value = -2379.125

$ws.Cells.Item(83, 8).Value = 422.36365
$ws.Cells.Item(83, 9).Value = 527
$ws.Cells.Item(83, 10).Value = 383.125
$ws.Cells.Item(83, 11).Value = 2635
$ws.Cells.Item(83, 12).Value = 1915.625
$ws.Cells.Item(83, 13).Value = 2357
$ws.Cells.Item(83, 14).Value = -11899.625

$ws.Cells.Item(134, 8).Value = 2915.675
$ws.Cells.Item(134, 9).Value = 951.8570999999999
$ws.Cells.Item(134, 10).Value = 7497.9165
$ws.Cells.Item(134, 11).Value = 2855.5713
$ws.Cells.Item(134, 12).Value = 22493.7495
$ws.Cells.Item(134, 13).Value = -320.5712999999996
$ws.Cells.Item(134, 14).Value = -27563.7495

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1053.6
$ws.Cells.Item(22, 9).Value = 785
$ws.Cells.Item(22, 10).Value = 1232.6666
$ws.Cells.Item(22, 11).Value = 785
$ws.Cells.Item(22, 12).Value = 1232.6666
$ws.Cells.Item(22, 13).Value = -435
$ws.Cells.Item(22, 14).Value = -1932.6666

$ws.Cells.Item(31, 8).Value = 2501.69
$ws.Cells.Item(31, 9).Value = 2349.0576
$ws.Cells.Item(31, 10).Value = 2667.0417
$ws.Cells.Item(31, 11).Value = 2349.0576
$ws.Cells.Item(31, 12).Value = 2667.0417
$ws.Cells.Item(31, 13).Value = -2054.0576
$ws.Cells.Item(31, 14).Value = -3257.0417

$ws.Cells.Item(34, 8).Value = 2501.69
$ws.Cells.Item(34, 9).Value = 2349.0576
$ws.Cells.Item(34, 10).Value = 2667.0417
$ws.Cells.Item(34, 11).Value = 2349.0576
$ws.Cells.Item(34, 12).Value = 2667.0417
$ws.Cells.Item(34, 13).Value = -2147.0576
$ws.Cells.Item(34, 14).Value = -3071.0417

$ws.Cells.Item(58, 8).Value = 910926.4399999999
$ws.Cells.Item(58, 9).Value = 1251399
$ws.Cells.Item(58, 10).Value = 2999.6667
$ws.Cells.Item(58, 11).Value = 1251399
$ws.Cells.Item(58, 12).Value = 2999.6667
$ws.Cells.Item(58, 13).Value = -1251196
$ws.Cells.Item(58, 14).Value = -3405.6667

$ws.Cells.Item(100, 8).Value = 55000
$ws.Cells.Item(100, 10).Value = 55000
$ws.Cells.Item(100, 12).Value = 55000
$ws.Cells.Item(100, 14).Value = -57164

$ws.Cells.Item(105, 8).Value = 1423096.8
$ws.Cells.Item(105, 9).Value = 2843003.2
$ws.Cells.Item(105, 11).Value = 2843003.2
$ws.Cells.Item(105, 13).Value = -2841256.2

$ws.Cells.Item(132, 8).Value = 15164762
$ws.Cells.Item(132, 9).Value = 18523552
$ws.Cells.Item(132, 11).Value = 55570656
$ws.Cells.Item(132, 13).Value = -55568126

$ws.Cells.Item(134, 8).Value = 2678.5881
$ws.Cells.Item(134, 9).Value = 2572.2666
$ws.Cells.Item(134, 10).Value = 3476
$ws.Cells.Item(134, 11).Value = 7716.7998
$ws.Cells.Item(134, 12).Value = 10428
$ws.Cells.Item(134, 13).Value = -5181.7998
$ws.Cells.Item(134, 14).Value = -15498

$ws.Cells.Item(136, 8).Value = 910926.4399999999
$ws.Cells.Item(136, 9).Value = 1251399
$ws.Cells.Item(136, 10).Value = 2999.6667
$ws.Cells.Item(136, 11).Value = 3754197
$ws.Cells.Item(136, 12).Value = 8999.000100000001
$ws.Cells.Item(136, 13).Value = -3751647
$ws.Cells.Item(136, 14).Value = -14099.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 5842.0415
$ws.Cells.Item(56, 9).Value = 5842.0415
$ws.Cells.Item(56, 11).Value = 5842.0415
$ws.Cells.Item(56, 13).Value = -5312.0415

$ws.Cells.Item(68, 8).Value = 459712.06
$ws.Cells.Item(68, 10).Value = 557779.7
$ws.Cells.Item(68, 12).Value = 1673339.1
$ws.Cells.Item(68, 14).Value = -1674961.1

$ws.Cells.Item(71, 8).Value = 459712.06
$ws.Cells.Item(71, 10).Value = 557779.7
$ws.Cells.Item(71, 12).Value = 5020017.3
$ws.Cells.Item(71, 14).Value = -5028129.3

$ws.Cells.Item(107, 8).Value = 1999
$ws.Cells.Item(107, 9).Value = 1999
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 5997
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -4077
$ws.Cells.Item(107, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 413.4
$ws.Cells.Item(122, 9).Value = 389
$ws.Cells.Item(122, 10).Value = 423.85715
$ws.Cells.Item(122, 11).Value = 3501
$ws.Cells.Item(122, 12).Value = 3814.71435
$ws.Cells.Item(122, 13).Value = -1051
$ws.Cells.Item(122, 14).Value = -8714.71435

$ws.Cells.Item(137, 8).Value = 5798100.5
$ws.Cells.Item(137, 9).Value = 2731.3333
$ws.Cells.Item(137, 10).Value = 11013933
$ws.Cells.Item(137, 11).Value = 8193.999899999999
$ws.Cells.Item(137, 12).Value = 33041799
$ws.Cells.Item(137, 13).Value = -3093.999899999999
$ws.Cells.Item(137, 14).Value = -33051999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 9526301
$ws.Cells.Item(113, 9).Value = 15153502
$ws.Cells.Item(113, 10).Value = 3345
$ws.Cells.Item(113, 11).Value = 15153502
$ws.Cells.Item(113, 12).Value = 3345
$ws.Cells.Item(113, 13).Value = -15151332
$ws.Cells.Item(113, 14).Value = -7685

$ws.Cells.Item(126, 8).Value = 4097.45
$ws.Cells.Item(126, 9).Value = 1870
$ws.Cells.Item(126, 10).Value = 6324.9
$ws.Cells.Item(126, 11).Value = 5610
$ws.Cells.Item(126, 12).Value = 18974.7
$ws.Cells.Item(126, 13).Value = -3140
$ws.Cells.Item(126, 14).Value = -23914.7

$ws.Cells.Item(132, 8).Value = 4632.75
$ws.Cells.Item(132, 9).Value = 4062.5715
$ws.Cells.Item(132, 10).Value = 5963.1665
$ws.Cells.Item(132, 11).Value = 12187.7145
$ws.Cells.Item(132, 12).Value = 17889.4995
$ws.Cells.Item(132, 13).Value = -9657.7145
$ws.Cells.Item(132, 14).Value = -22949.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4345.59
$ws.Cells.Item(132, 9).Value = 4329.4346
$ws.Cells.Item(132, 10).Value = 5945
$ws.Cells.Item(132, 11).Value = 12988.3038
$ws.Cells.Item(132, 12).Value = 17835
$ws.Cells.Item(132, 13).Value = -10458.3038
$ws.Cells.Item(132, 14).Value = -22895

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 1000000000
$ws.Cells.Item(70, 9).Value = 1000000000
$ws.Cells.Item(70, 11).Value = 1000000000
$ws.Cells.Item(70, 13).Value = -999999685

$ws.Cells.Item(73, 8).Value = 1000000000
$ws.Cells.Item(73, 9).Value = 1000000000
$ws.Cells.Item(73, 11).Value = 1000000000
$ws.Cells.Item(73, 13).Value = -999998908

$ws.Cells.Item(126, 8).Value = 2141.8572
$ws.Cells.Item(126, 9).Value = 1123.25
$ws.Cells.Item(126, 10).Value = 3500
$ws.Cells.Item(126, 11).Value = 3369.75
$ws.Cells.Item(126, 12).Value = 10500
$ws.Cells.Item(126, 13).Value = -899.75
$ws.Cells.Item(126, 14).Value = -15440

$ws.Cells.Item(132, 8).Value = 16710761
$ws.Cells.Item(132, 9).Value = 50589.117
$ws.Cells.Item(132, 10).Value = 125001870
$ws.Cells.Item(132, 11).Value = 151767.351
$ws.Cells.Item(132, 12).Value = 375005610
$ws.Cells.Item(132, 13).Value = -149237.351
$ws.Cells.Item(132, 14).Value = -375010670

$ws.Cells.Item(136, 8).Value = 9531.08
$ws.Cells.Item(136, 9).Value = 6347.5454
$ws.Cells.Item(136, 10).Value = 9924.550999999999
$ws.Cells.Item(136, 11).Value = 19042.6362
$ws.Cells.Item(136, 12).Value = 29773.653
$ws.Cells.Item(136, 13).Value = -16492.6362
$ws.Cells.Item(136, 14).Value = -34873.653
